$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 15.55806759907682
$ws.Range("C2").Value = 8.540528589571126
$ws.Range("E2").Value = 14.10198856038042
$ws.Range("F2").Value = 45.7121609316297
$ws.Range("G2").Value = 51.62137166458338
$ws.Range("H2").Value = 19.9799272873642
$ws.Range("J2").Value = 9.608903133348065
$ws.Range("K2").Value = 11.41440754542955
$ws.Range("L2").Value = 11.47510604193225
# Row 3
$ws.Range("B3").Value = 15.4060128603375
$ws.Range("C3").Value = 8.522617624406783
$ws.Range("E3").Value = 14.10863788487438
$ws.Range("F3").Value = 45.67171881958199
$ws.Range("G3").Value = 51.57878513774838
$ws.Range("H3").Value = 20.01264539227814
$ws.Range("J3").Value = 9.620583546541793
$ws.Range("K3").Value = 11.31511672461193
$ws.Range("L3").Value = 11.47221032041698
# Row 4
$ws.Range("B4").Value = 15.3157340377111
$ws.Range("C4").Value = 8.51130705454386
$ws.Range("E4").Value = 14.114740175127
$ws.Range("F4").Value = 45.65616804600596
$ws.Range("G4").Value = 51.5640173244451
$ws.Range("H4").Value = 20.03566352009329
$ws.Range("J4").Value = 9.628055132384107
$ws.Range("K4").Value = 11.25642730046902
$ws.Range("L4").Value = 11.47221289140771
# Row 5
$ws.Range("B5").Value = 15.2797606330304
$ws.Range("C5").Value = 8.506619269600597
$ws.Range("E5").Value = 14.11773579936433
$ws.Range("F5").Value = 45.65216898181959
$ws.Range("G5").Value = 51.56086162301161
$ws.Range("H5").Value = 20.04577933531914
$ws.Range("J5").Value = 9.631175510347667
$ws.Range("K5").Value = 11.23310809404046
$ws.Range("L5").Value = 11.47266319672898
# Row 6
$ws.Range("B6").Value = 15.27383771305151
$ws.Range("C6").Value = 8.505836120344805
$ws.Range("E6").Value = 14.11826398455059
$ws.Range("F6").Value = 45.65164624418327
$ws.Range("G6").Value = 51.56051047496155
$ws.Range("H6").Value = 20.04750347065896
$ws.Range("J6").Value = 9.631698223422106
$ws.Range("K6").Value = 11.22927273651169
$ws.Range("L6").Value = 11.47276514471121
# Row 7
$ws.Range("B7").Value = 15.3152455315993
$ws.Range("C7").Value = 8.511244151673433
$ws.Range("E7").Value = 14.1147785133766
$ws.Range("F7").Value = 45.65610464232586
$ws.Range("G7").Value = 51.56396317702343
$ws.Range("H7").Value = 20.03579696768388
$ws.Range("J7").Value = 9.628096908255841
$ws.Range("K7").Value = 11.25611035921157
$ws.Range("L7").Value = 11.47221714347032
# Row 8
$ws.Range("B8").Value = 15.50502399334272
$ws.Range("C8").Value = 8.53441794997361
$ws.Range("E8").Value = 14.103862749798
$ws.Range("F8").Value = 45.69629324829801
$ws.Range("G8").Value = 51.60432674152512
$ws.Range("H8").Value = 19.99060010148719
$ws.Range("J8").Value = 9.612868485898758
$ws.Range("K8").Value = 11.37971511432415
$ws.Range("L8").Value = 11.47373900509994
# Row 9
$ws.Range("B9").Value = 15.89979042705104
$ws.Range("C9").Value = 8.577373314461232
$ws.Range("E9").Value = 14.09843146323078
$ws.Range("F9").Value = 45.84849447373166
$ws.Range("G9").Value = 51.77367260442076
$ws.Range("H9").Value = 19.92524390571324
$ws.Range("J9").Value = 9.585371568249958
$ws.Range("K9").Value = 11.63901075704613
$ws.Range("L9").Value = 11.49077837001723
# Row 10
$ws.Range("B10").Value = 16.20104904796498
$ws.Range("C10").Value = 8.607412881171776
$ws.Range("E10").Value = 14.1041069303327
$ws.Range("F10").Value = 46.00463687396758
$ws.Range("G10").Value = 51.95281507287211
$ws.Range("H10").Value = 19.89145885927578
$ws.Range("J10").Value = 9.566593963042495
$ws.Range("K10").Value = 11.83822591510056
$ws.Range("L10").Value = 11.51175880344585
# Row 11
$ws.Range("B11").Value = 16.33997791702568
$ws.Range("C11").Value = 8.620748177931683
$ws.Range("E11").Value = 14.10876952860266
$ws.Range("F11").Value = 46.08517162678311
$ws.Range("G11").Value = 52.04609593925527
$ws.Range("H11").Value = 19.87918703094084
$ws.Range("J11").Value = 9.558356966630315
$ws.Range("K11").Value = 11.93039372196398
$ws.Range("L11").Value = 11.52311431806778
# Row 12
$ws.Range("B12").Value = 16.39281057014234
$ws.Range("C12").Value = 8.625750379557248
$ws.Range("E12").Value = 14.11083254561269
$ws.Range("F12").Value = 46.117021739446
$ws.Range("G12").Value = 52.08310206441757
$ws.Range("H12").Value = 19.87498579252819
$ws.Range("J12").Value = 9.555281405994874
$ws.Range("K12").Value = 11.96548730967629
$ws.Range("L12").Value = 11.52767226862591
# Row 13
$ws.Range("B13").Value = 16.38142293907071
$ws.Range("C13").Value = 8.624675184772013
$ws.Range("E13").Value = 14.1103750419754
$ws.Range("F13").Value = 46.11010228769168
$ws.Range("G13").Value = 52.07505751556467
$ws.Range("H13").Value = 19.87587076937241
$ws.Range("J13").Value = 9.5559418473871
$ws.Range("K13").Value = 11.95792123839039
$ws.Range("L13").Value = 11.52667920746196
# Row 14
$ws.Range("B14").Value = 16.34432028846445
$ws.Range("C14").Value = 8.62116066663819
$ws.Range("E14").Value = 14.10893330503214
$ws.Range("F14").Value = 46.08776491406869
$ws.Range("G14").Value = 52.04910681058074
$ws.Range("H14").Value = 19.87883245500183
$ws.Range("J14").Value = 9.558103065813084
$ws.Range("K14").Value = 11.93327722590732
$ws.Range("L14").Value = 11.52348415106368
# Row 15
$ws.Range("B15").Value = 16.3216215114417
$ws.Range("C15").Value = 8.619001721172086
$ws.Range("E15").Value = 14.10808887216765
$ws.Range("F15").Value = 46.07425844849254
$ws.Range("G15").Value = 52.03342998504807
$ws.Range("H15").Value = 19.88070464610941
$ws.Range("J15").Value = 9.559432546502956
$ws.Range("K15").Value = 11.91820611164527
$ws.Range("L15").Value = 11.52156058747128
# Row 16
$ws.Range("B16").Value = 16.19200356341507
$ws.Range("C16").Value = 8.606534737568023
$ws.Range("E16").Value = 14.10384393694409
$ws.Range("F16").Value = 45.99956388729429
$ws.Range("G16").Value = 51.94695514113854
$ws.Range("H16").Value = 19.89232321744173
$ws.Range("J16").Value = 9.567138386566048
$ws.Range("K16").Value = 11.83223100462022
$ws.Range("L16").Value = 11.51105292474399
# Row 17
$ws.Range("B17").Value = 16.11293589913337
$ws.Range("C17").Value = 8.598802054508747
$ws.Range("E17").Value = 14.1017714973864
$ws.Range("F17").Value = 45.95616654605185
$ws.Range("G17").Value = 51.89691666452708
$ws.Range("H17").Value = 19.90024445945287
$ws.Range("J17").Value = 9.571943615916407
$ws.Range("K17").Value = 11.77986190581896
$ws.Range("L17").Value = 11.50506895776019
# Row 18
$ws.Range("B18").Value = 16.06763853321083
$ws.Range("C18").Value = 8.594323381634029
$ws.Range("E18").Value = 14.1007754580257
$ws.Range("F18").Value = 45.93210108354167
$ws.Range("G18").Value = 51.86924628278228
$ws.Range("H18").Value = 19.90509204622687
$ws.Range("J18").Value = 9.574736182612995
$ws.Range("K18").Value = 11.74988772652466
$ws.Range("L18").Value = 11.50179780497947
# Row 19
$ws.Range("B19").Value = 16.05233404033649
$ws.Range("C19").Value = 8.592801653795854
$ws.Range("E19").Value = 14.10047192798333
$ws.Range("F19").Value = 45.92410712473188
$ws.Range("G19").Value = 51.8600686277177
$ws.Range("H19").Value = 19.90678340483287
$ws.Range("J19").Value = 9.575686639108627
$ws.Range("K19").Value = 11.73976517959437
$ws.Range("L19").Value = 11.50071963768311
# Row 20
$ws.Range("B20").Value = 16.12133447664834
$ws.Range("C20").Value = 8.599628425268593
$ws.Range("E20").Value = 14.10197184240219
$ws.Range("F20").Value = 45.96069367426128
$ws.Range("G20").Value = 51.90212850398819
$ws.Range("H20").Value = 19.899371057003
$ws.Range("J20").Value = 9.571429119850512
$ws.Range("K20").Value = 11.78542167689102
$ws.Range("L20").Value = 11.50568831625269
# Row 21
$ws.Range("B21").Value = 16.35521255522369
$ws.Range("C21").Value = 8.622194258495911
$ws.Range("E21").Value = 14.10934872223613
$ws.Range("F21").Value = 46.09428932743338
$ws.Range("G21").Value = 52.05668360354585
$ws.Range("H21").Value = 19.87795043306387
$ws.Range("J21").Value = 9.557467081980571
$ws.Range("K21").Value = 11.94051081094973
$ws.Range("L21").Value = 11.52441563998456
# Row 22
$ws.Range("B22").Value = 16.50934506017941
$ws.Range("C22").Value = 8.636664897310151
$ws.Range("E22").Value = 14.11590242101664
$ws.Range("F22").Value = 46.18948283874871
$ws.Range("G22").Value = 52.16749503100429
$ws.Range("H22").Value = 19.86654963411249
$ws.Range("J22").Value = 9.548596166091881
$ws.Range("K22").Value = 12.04297346157998
$ws.Range("L22").Value = 11.53815671296588
# Row 23
$ws.Range("B23").Value = 16.42698025392845
$ws.Range("C23").Value = 8.628967081206511
$ws.Range("E23").Value = 14.11224670644883
$ws.Range("F23").Value = 46.13795999520185
$ws.Range("G23").Value = 52.10746076626403
$ws.Range("H23").Value = 19.8723965366056
$ws.Range("J23").Value = 9.553307575192129
$ws.Range("K23").Value = 11.98819633443004
$ws.Range("L23").Value = 11.53068635074002
# Row 24
$ws.Range("B24").Value = 16.11753697869233
$ws.Range("C24").Value = 8.599254926003534
$ws.Range("E24").Value = 14.10188065754208
$ws.Range("F24").Value = 45.95864420414809
$ws.Range("G24").Value = 51.89976881196924
$ws.Range("H24").Value = 19.89976500780681
$ws.Range("J24").Value = 9.571661630116406
$ws.Range("K24").Value = 11.78290768618373
$ws.Range("L24").Value = 11.50540777733031
# Row 25
$ws.Range("B25").Value = 15.79083611525298
$ws.Range("C25").Value = 8.566019430156976
$ws.Range("E25").Value = 14.09819754024588
$ws.Range("F25").Value = 45.79950010621405
$ws.Range("G25").Value = 51.7182227528285
$ws.Range("H25").Value = 19.94042792154344
$ws.Range("J25").Value = 9.592558781969799
$ws.Range("K25").Value = 11.56721656610206
$ws.Range("L25").Value = 11.47266319672898
